$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BP-813: Affiliate Mapping for True Independent Stations
# Rename header columns:
#   H1: "Affiliation Mismatch Note" -> "IsTrueIND"
#   J1: "SalesGroupName" -> "RepFirm"
$ws.Range("H1").Value = "IsTrueIND"
$ws.Range("J1").Value = "RepFirm"

# Extend the bordered header/body formatting (already used by column H)
# across the I:J columns to match the rest of the table.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("I2:J5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
